# Update TPM-derived NATMI metrics in the Dlk1-Notch4 LR-pairs sheet.
# The underlying TPM values for the receptor (Notch4) in the "ECs" target
# cluster were refreshed, which cascades through the specificity and edge
# weight columns (O,P,Q,R,S,T) for every row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending=FAPs, Target=ECs)
$ws.Range("M2").Value = 31.618405
$ws.Range("N2").Value = 94.855215
$ws.Range("O2").Value = 0.8578613706944929
$ws.Range("P2").Value = 0.8578613706944929
$ws.Range("Q2").Value = 133.5730691061433
$ws.Range("R2").Value = 1202.15762195529
$ws.Range("S2").Value = 0.6042111069796038
$ws.Range("T2").Value = 0.6042111069796039

# Row 3 (Sending=FAPs, Target=FAPs)
$ws.Range("O3").Value = 0.08747555172986397
$ws.Range("P3").Value = 0.08747555172986396
$ws.Range("S3").Value = 0.06161100353727818
$ws.Range("T3").Value = 0.06161100353727818

# Row 4 (Sending=FAPs, Target=MuSCs)
$ws.Range("M4").Value = 2.014730333333334
$ws.Range("N4").Value = 6.044191000000001
$ws.Range("O4").Value = 0.05466307757564324
$ws.Range("P4").Value = 0.05466307757564324
$ws.Range("Q4").Value = 8.511299480305111
$ws.Range("R4").Value = 76.60169532274601
$ws.Range("S4").Value = 0.03850043811408955
$ws.Range("T4").Value = 0.03850043811408956

# Row 5 (Sending=MuSCs, Target=ECs)
$ws.Range("M5").Value = 31.618405
$ws.Range("N5").Value = 94.855215
$ws.Range("O5").Value = 0.8578613706944929
$ws.Range("P5").Value = 0.8578613706944929
$ws.Range("Q5").Value = 56.07451404418499
$ws.Range("R5").Value = 504.6706263976649
$ws.Range("S5").Value = 0.2536502637148891
$ws.Range("T5").Value = 0.2536502637148891

# Row 6 (Sending=MuSCs, Target=FAPs)
$ws.Range("O6").Value = 0.08747555172986397
$ws.Range("P6").Value = 0.08747555172986396
$ws.Range("S6").Value = 0.02586454819258579
$ws.Range("T6").Value = 0.02586454819258579

# Row 7 (Sending=MuSCs, Target=MuSCs)
$ws.Range("M7").Value = 2.014730333333334
$ws.Range("N7").Value = 6.044191000000001
$ws.Range("O7").Value = 0.05466307757564324
$ws.Range("P7").Value = 0.05466307757564324
$ws.Range("Q7").Value = 3.573077907369
$ws.Range("R7").Value = 32.157701166321
$ws.Range("S7").Value = 0.01616263946155369
$ws.Range("T7").Value = 0.01616263946155369

$wb.Save()
